$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.183.99"
$ws.Range("E2").Value = "  +3.01%  "
$ws.Range("D3").Value = "3.589.66"
$ws.Range("E3").Value = "  +2.84%  "
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "625.43"
$ws.Range("E5").Value = "  +2.75%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "157.89"
$ws.Range("E6").Value = "  +6.35%  "
$ws.Range("D7").Value = "3.582.64"
$ws.Range("E7").Value = "  +2.64%  "
$ws.Range("E8").Value = "  +0.16%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.494"
$ws.Range("E9").Value = "  +2.60%  "
$ws.Range("E10").Value = "  +8.23%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "7.39"
$ws.Range("E11").Value = "  +7.01%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.441"
$ws.Range("E12").Value = "  +4.60%  "
$ws.Range("E13").Value = "  +5.18%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "33.61"
$ws.Range("E14").Value = "  +7.11%  "
$ws.Range("D15").Value = "4.199.84"
$ws.Range("E15").Value = "  +2.90%  "
$ws.Range("D16").Value = "69.353.40"
$ws.Range("E16").Value = "  +3.40%  "
$ws.Range("D17").Value = "3.582.20"
$ws.Range("E17").Value = "  +2.88%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.81"
$ws.Range("E19").Value = "  +5.79%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "16.16"
$ws.Range("E20").Value = "  +7.51%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.17"
$ws.Range("E21").Value = "  +12.63%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "461.96"
$ws.Range("E22").Value = "  +4.11%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.646"
$ws.Range("E23").Value = "  +3.54%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "78.90"
$ws.Range("E24").Value = "  +2.31%  "
$ws.Range("E25").Value = "  +6.59%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "10.69"
$ws.Range("E26").Value = "  +5.49%  "
$ws.Range("D27").Value = "3.736.33"
$ws.Range("E27").Value = "  +2.93%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.00"
$ws.Range("E28").Value = "  -0.11%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.25"
$ws.Range("E29").Value = "  +11.62%  "
$ws.Range("E30").Value = "  +4.06%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.72"
$ws.Range("E31").Value = "  +9.31%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.173"
$ws.Range("E32").Value = "  +4.87%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.999"
$ws.Range("E33").Value = "  -0.15%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.52"
$ws.Range("E34").Value = "  +6.52%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "26.44"
$ws.Range("E35").Value = "  +3.19%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.94"
$ws.Range("E36").Value = "  +4.72%  "
$ws.Range("D37").Value = "3.588.07"
$ws.Range("E37").Value = "  +3.28%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "8.40"
$ws.Range("E38").Value = "  +5.34%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.41"
$ws.Range("E39").Value = "  +9.62%  "
$ws.Range("E40").Value = "  +0.00%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "179.89"
$ws.Range("E41").Value = "  +6.32%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0925"
$ws.Range("E42").Value = "  +6.48%  "
$ws.Range("E43").Value = "  +0.01%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.68"
$ws.Range("E44").Value = "  +4.53%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "31.28"
$ws.Range("E45").Value = "  +19.26%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.907"
$ws.Range("E46").Value = "  +2.92%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.38"
$ws.Range("E47").Value = "  +11.14%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "45.99"
$ws.Range("E48").Value = "  +0.77%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.75"
$ws.Range("E49").Value = "  +9.16%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.83"
$ws.Range("E50").Value = "  +3.67%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.267"
$ws.Range("E51").Value = "  +8.53%  "
